$wb = $excel.ActiveWorkbook

# Sheet: Manila Philippines
$ws = $wb.Worksheets.Item("Manila Philippines")
$ws.Range("E2").Value = 0.3333
$ws.Range("E3").Value = 0.3333
$ws.Range("E4").Value = 0.3333
$ws.Range("O4:W4").Value = 0

# Sheet: Milwaukee Pmc Hq Wisconsin
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("O7").ClearContents()

# Sheet: Milwaukee Wisconsin
$ws = $wb.Worksheets.Item("Milwaukee Wisconsin")
$ws.Range("O5").ClearContents()

# Sheet: South Beloit Gardner St Illinois
$ws = $wb.Worksheets.Item("South Beloit Gardner St Illino")
$ws.Range("O7").ClearContents()

# Sheet: Rock Road Radford Virginia
$ws = $wb.Worksheets.Item("Rock Road Radford Virginia")
$ws.Range("O2").ClearContents()
$ws.Range("O3").ClearContents()
